$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "26.914.79"
$ws.Cells.Item(2, 5).Value = "  -0.13%  "

$ws.Cells.Item(3, 4).Value = "1.549.93"
$ws.Cells.Item(3, 5).Value = "  -0.44%  "

Set-TextValue $ws.Cells.Item(5, 4) "206.30"

$ws.Cells.Item(6, 5).Value = "  +0.11%  "

Set-TextValue $ws.Cells.Item(8, 4) "22.14"
$ws.Cells.Item(8, 5).Value = "  +2.87%  "

$ws.Cells.Item(9, 5).Value = "  -0.51%  "

$ws.Cells.Item(10, 5).Value = "  +0.62%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.0856"
$ws.Cells.Item(11, 5).Value = "  -0.30%  "

$ws.Cells.Item(12, 4).Value = "1.771.22"
$ws.Cells.Item(12, 5).Value = "  -0.43%  "

$ws.Cells.Item(13, 4).Value = "1.545.40"
$ws.Cells.Item(13, 5).Value = "  -0.71%  "

$ws.Cells.Item(14, 5).Value = "  +0.58%  "

Set-TextValue $ws.Cells.Item(15, 4) "0.517"
$ws.Cells.Item(15, 5).Value = "  +0.45%  "

$ws.Cells.Item(16, 4).Value = "26.905.56"
$ws.Cells.Item(16, 5).Value = "  -0.13%  "

Set-TextValue $ws.Cells.Item(17, 4) "61.61"
$ws.Cells.Item(17, 5).Value = "  -0.18%  "

Set-TextValue $ws.Cells.Item(18, 4) "217.17"
$ws.Cells.Item(18, 5).Value = "  +1.21%  "

$ws.Cells.Item(19, 4).Value = "0.0₃0697"
$ws.Cells.Item(19, 5).Value = "  +1.36%  "

Set-TextValue $ws.Cells.Item(20, 4) "7.25"
$ws.Cells.Item(20, 5).Value = "  -0.13%  "

Set-TextValue $ws.Cells.Item(22, 4) "4.05"
$ws.Cells.Item(22, 5).Value = "  +0.00%  "

Set-TextValue $ws.Cells.Item(23, 4) "9.23"
$ws.Cells.Item(23, 5).Value = "  +0.36%  "

$ws.Cells.Item(24, 5).Value = "  -1.08%  "

Set-TextValue $ws.Cells.Item(25, 4) "153.96"
$ws.Cells.Item(25, 5).Value = "  +0.36%  "

Set-TextValue $ws.Cells.Item(26, 4) "6.61"
$ws.Cells.Item(26, 5).Value = "  -0.86%  "

Set-TextValue $ws.Cells.Item(27, 4) "14.92"
$ws.Cells.Item(27, 5).Value = "  +0.22%  "

$ws.Cells.Item(28, 5).Value = "  +0.61%  "

$ws.Cells.Item(29, 5).Value = "  -0.28%  "

$ws.Cells.Item(30, 5).Value = "  +1.70%  "

Set-TextValue $ws.Cells.Item(31, 4) "1.09"
$ws.Cells.Item(31, 5).Value = "  -0.69%  "

$ws.Cells.Item(32, 5).Value = "  -0.71%  "

$ws.Cells.Item(33, 4).Value = "1.421.20"
$ws.Cells.Item(33, 5).Value = "  +3.73%  "

Set-TextValue $ws.Cells.Item(34, 4) "3.08"
$ws.Cells.Item(34, 5).Value = "  +3.99%  "

$ws.Cells.Item(35, 5).Value = "  +2.15%  "

$ws.Cells.Item(36, 5).Value = "  -0.29%  "

$ws.Cells.Item(37, 5).Value = "  +0.08%  "

$ws.Cells.Item(38, 5).Value = "  +0.50%  "

Set-TextValue $ws.Cells.Item(39, 4) "0.523"
$ws.Cells.Item(39, 5).Value = "  +0.19%  "

Set-TextValue $ws.Cells.Item(40, 4) "0.807"
$ws.Cells.Item(40, 5).Value = "  -0.21%  "

$ws.Cells.Item(41, 5).Value = "  -0.38%  "

$ws.Cells.Item(42, 5).Value = "  +3.43%  "

$ws.Cells.Item(43, 5).Value = "  +3.11%  "

$ws.Cells.Item(44, 5).Value = "  +0.83%  "

Set-TextValue $ws.Cells.Item(45, 4) "64.44"
$ws.Cells.Item(45, 5).Value = "  +1.12%  "

$ws.Cells.Item(46, 5).Value = "  +1.08%  "

$ws.Cells.Item(47, 4).Value = "1.684.92"
$ws.Cells.Item(47, 5).Value = "  -0.42%  "

Set-TextValue $ws.Cells.Item(48, 4) "87.41"
$ws.Cells.Item(48, 5).Value = "  +1.56%  "

$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(49, 4) "0.0515"
$ws.Cells.Item(49, 5).Value = "  +1.40%  "

$ws.Cells.Item(50, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(50, 4).Value = "0.0₆0101"
$ws.Cells.Item(50, 5).Value = "  +3.74%  "

Set-TextValue $ws.Cells.Item(51, 4) "0.0959"
$ws.Cells.Item(51, 5).Value = "  +0.39%  "
